# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (F column) counts for several rows across the
# "展览", "演出" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9827
$ws1.Range("F13").Value = 3042
$ws1.Range("F14").Value = 2298
$ws1.Range("F16").Value = 1995
$ws1.Range("F22").Value = 32
$ws1.Range("F25").Value = 39
$ws1.Range("F33").Value = 232
$ws1.Range("F34").Value = 1568
$ws1.Range("F37").Value = 37
$ws1.Range("F38").Value = 405
$ws1.Range("F39").Value = 839
$ws1.Range("F41").Value = 328

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9827
$ws4.Range("F3").Value = 417
$ws4.Range("F7").Value = 1
$ws4.Range("F15").Value = 3042
$ws4.Range("F16").Value = 2298
$ws4.Range("F18").Value = 1995
$ws4.Range("F24").Value = 32
$ws4.Range("F27").Value = 39
$ws4.Range("F39").Value = 232
$ws4.Range("F40").Value = 1568
$ws4.Range("F44").Value = 37
$ws4.Range("F45").Value = 405
$ws4.Range("F46").Value = 839
$ws4.Range("F48").Value = 328
